# Actualización automática 2025-10-15 08:30:09
#
# Updates sales figures for advisor "LINDAO ZUÑIGA BRYAN JOSE" on three
# related sheets: the detail-by-group sheet, the monthly-sales sheet, and
# the monthly-compliance summary sheet. None of the cells are live Excel
# formulas (the workbook stores pre-computed static values), so every
# downstream total / subtotal / percentage that depends on the changed
# leaf values is recomputed here and written explicitly as well.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" (sales by product group, per client row)
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Client "BORJA TORRES LETTY JANET" (row 11): FREGADEROS DE COCINA (col E)
$wsGrupo.Range("E11").Value = 1591.79

# Client "CASA FERRETERIA FONG S.A. FERREFONG" (row 16): NO RESURTIBLES (col P)
$wsGrupo.Range("P16").Value = 1445.41

# Client "MUÑOZ FALCONES SERGIO BACILIO" (row 41): INODOROS (col H)
$wsGrupo.Range("H41").Value = 71.09999999999999

# Row 60 holds "<n> de 58" counts of non-zero entries per column; update the
# three columns whose non-zero-cell count changed because of the edits above.
$wsGrupo.Range("E60").Value = "1 de 58"
$wsGrupo.Range("H60").Value = "2 de 58"
$wsGrupo.Range("P60").Value = "1 de 58"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" (monthly sales, per client row)
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# "octubre" (October) column (F) for the same three clients
$wsMensual.Range("F11").Value = 1591.79
$wsMensual.Range("F16").Value = 1445.41
$wsMensual.Range("F41").Value = 71.09999999999999

# Row 60 totals row: October total increases by the sum of the three deltas
$wsMensual.Range("F60").Value = 16201.25

# ---------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL" (monthly compliance by product group)
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen the "POR CUMPLIR" (column E) to fit the new, larger values.
# (ColumnWidth and the stored OOXML column width differ by a constant
# 5/6-character padding offset in this engine, so back that out here to
# land on exactly width="24" in the saved file.)
$wsCumpl.Columns("E").ColumnWidth = 23.166666666666668

# Row 4: FREGADEROS DE COCINA — VENTA (D), POR CUMPLIR (E), CUMPLIMIENTO (F)
$wsCumpl.Range("D4").Value = 1448.05
$wsCumpl.Range("E4").Value = -144.0213934184001
$wsCumpl.Range("F4").Value = 1.110443430988788

# Row 6: INODOROS
$wsCumpl.Range("D6").Value = 465.47
$wsCumpl.Range("E6").Value = 384.37419682004
$wsCumpl.Range("F6").Value = 0.5477121591718844

# Row 8: NO RESURTIBLES
$wsCumpl.Range("D8").Value = 1445.41
$wsCumpl.Range("E8").Value = -929.2881264521661
$wsCumpl.Range("F8").Value = 2.800520718225363

# Row 14: TOTAL
$wsCumpl.Range("D14").Value = 16201.25
$wsCumpl.Range("E14").Value = 38282.51774946896
$wsCumpl.Range("F14").Value = 0.2973592075074123
